$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Root_surface_diff" -> "Root_proj_area_diff" and
# "Root_surface_avg" -> "Root_proj_area_avg" across the whole sheet
# (commit message: "update all figures from root surf to root area")
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value()
    if ($val -eq "Root_surface_diff") {
        $cell.Value = "Root_proj_area_diff"
    } elseif ($val -eq "Root_surface_avg") {
        $cell.Value = "Root_proj_area_avg"
    }
}

# Restore the active selection seen in the saved workbook
$ws.Range("D8").Select()
